# Form the consolidated report: populate the "Absent" (H) column.
# H is 1 when the "Real" (E) attendance count for that day is 0,
# and 0 otherwise. This also fills in the previously-blank H cells
# (rows where the day was fully "Real") with an explicit 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 3; $r -le $lastRow; $r++) {
    $real = $ws.Cells.Item($r, 5).Value2
    if ($null -eq $real) { $real = 0 }

    if ($real -eq 0) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
